$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3rd run")
$ws.Range("E21").Formula = '=(D21-C21)*$E$43'
